$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new rows (66-68) to the "pre-list of identified missing CMIP6
# requested variables" sheet: the new MODIS pextra table 126 variables
# (cltmodis, clwmodis, climodis) that must be deactivated by default and
# put on the identified-missing list (issue #746).

$ws.Range("A66").Value = "Emon"
$ws.Range("B66").Value = "cltmodis"
$ws.Range("C66").Value = 1
$ws.Range("H66").Value = "Implemented in EC-Earth table 126: grib code 47.126"

$ws.Range("A67").Value = "Emon"
$ws.Range("B67").Value = "clwmodis"
$ws.Range("C67").Value = 1
$ws.Range("H67").Value = "Implemented in EC-Earth table 126: grib code 48.126"

$ws.Range("A68").Value = "Emon"
$ws.Range("B68").Value = "climodis"
$ws.Range("C68").Value = 1
$ws.Range("H68").Value = "Implemented in EC-Earth table 126: grib code 49.126"

# Give the two latest comment cells the plain (non-charset) Calibri font
# that is introduced for these new entries.
$ws.Range("H67:H68").Font.Name = "Calibri"
$ws.Range("H67:H68").Font.Size = 11
$ws.Range("H67:H68").Font.Color = 0

# Update the view so the newly added rows are visible/selected, matching
# the author's saved selection state (rows 66-68 selected, scrolled so
# row 46 is at the top).
$ws.Rows("66:68").Select()
$excel.ActiveWindow.ScrollRow = 46
